# Add a new "image_leaf" trait row to the traits sheet, as part of moving
# all image-bearing views over to the new JSON based galleria data provider.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new row of data (name, type, description) right below the
# existing trait rows.
$ws.Range("A4").Value = "image_leaf"
$ws.Range("B4").Value = "image"
$ws.Range("C4").Value = "leaf image"

# Move the selection/active cell onto the newly added row.
[void]$ws.Range("A4").Select()
